$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 15 ("DDoS Attack: On a Web Server hosted on a Virtual Machine")
#    Rename the two DPS comparison labels.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(15)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -eq -1 -and $sh.TextFrame.HasText -eq -1) {
        $txt = $sh.TextFrame.TextRange.Text
        if ($txt -eq "DPS : Basic") {
            $sh.TextFrame.TextRange.Text = "With- out DPS"
        } elseif ($txt -eq "DPS : Standard") {
            $sh.TextFrame.TextRange.Text = "With DPS"
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Refresh the cached "datetimeFigureOut" footer field text (4/16/2018 ->
#    4/19/2018) on every slide layout and on the slide master.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapeRange) {
    for ($k = 1; $k -le $shapeRange.Count; $k++) {
        $shp = $shapeRange.Item($k)
        if ($shp.Type -eq 14) {
            $isDate = $false
            try { $isDate = ($shp.PlaceholderFormat.Type -eq 16) } catch {}
            if ($isDate -and $shp.HasTextFrame -eq -1 -and $shp.TextFrame.HasText -eq -1) {
                if ($shp.TextFrame.TextRange.Text -eq "4/16/2018") {
                    $shp.TextFrame.TextRange.Text = "4/19/2018"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $lo = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $lo.Shapes
}
